# Apply updated currentAveragePrice / LevePrice / LeveProfit figures
# (refreshed market-board data) across the ALC/ARM/BSM/CRP/CUL/GSM/LTW sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 4
$ws.Range("H4").Value = 2079.25
$ws.Range("I4").Value = 105.666664
$ws.Range("K4").Value = 105.666664
$ws.Range("M4").Value = 8.333336000000003
# Row 18
$ws.Range("H18").Value = 5321.304
$ws.Range("I18").Value = 5408.6816
$ws.Range("J18").Value = 3399
$ws.Range("K18").Value = 5408.6816
$ws.Range("L18").Value = 3399
$ws.Range("M18").Value = -5124.6816
$ws.Range("N18").Value = -3967
# Row 32
$ws.Range("H32").Value = 6049.625
$ws.Range("I32").Value = 4849.25
$ws.Range("J32").Value = 7250
$ws.Range("K32").Value = 4849.25
$ws.Range("L32").Value = 7250
$ws.Range("M32").Value = -4523.25
$ws.Range("N32").Value = -7902
# Row 39
$ws.Range("H39").Value = 101.545456
$ws.Range("I39").Value = 83.375
$ws.Range("K39").Value = 250.125
$ws.Range("M39").Value = 45.875
# Row 40
$ws.Range("H40").Value = 4453.5386
$ws.Range("I40").Value = 4314.143
$ws.Range("J40").Value = 4616.1665
$ws.Range("K40").Value = 4314.143
$ws.Range("L40").Value = 4616.1665
$ws.Range("M40").Value = -4139.143
$ws.Range("N40").Value = -4966.1665
# Row 43
$ws.Range("H43").Value = 5046.9
$ws.Range("I43").Value = 4069.1
$ws.Range("K43").Value = 4069.1
$ws.Range("M43").Value = -4000.1
# Row 45
$ws.Range("H45").Value = 2435.75
$ws.Range("I45").Value = 314
$ws.Range("J45").Value = 3143
$ws.Range("K45").Value = 942
$ws.Range("L45").Value = 9429
$ws.Range("M45").Value = -750
$ws.Range("N45").Value = -9813
# Row 46
$ws.Range("H46").Value = 128874.5
$ws.Range("J46").Value = 171349.33
$ws.Range("L46").Value = 514047.99
$ws.Range("N46").Value = -514285.99
# Row 47
$ws.Range("H47").Value = 7702
$ws.Range("I47").Value = 7702
$ws.Range("J47").Value = 0
$ws.Range("K47").Value = 7702
$ws.Range("L47").Value = 0
$ws.Range("M47").Value = -6730
$ws.Range("N47").Value = $null
# Row 49
$ws.Range("H49").Value = 4116.75
$ws.Range("J49").Value = 5039
$ws.Range("L49").Value = 15117
$ws.Range("N49").Value = -15389
# Row 54
$ws.Range("H54").Value = 10992.5
$ws.Range("J54").Value = 11391
$ws.Range("L54").Value = 11391
$ws.Range("N54").Value = -12363
# Row 60
$ws.Range("H60").Value = 128874.5
$ws.Range("J60").Value = 171349.33
$ws.Range("L60").Value = 514047.99
$ws.Range("N60").Value = -515015.99
# Row 64
$ws.Range("H64").Value = 45280
$ws.Range("J64").Value = 5452.25
$ws.Range("L64").Value = 5452.25
$ws.Range("N64").Value = -5948.25
# Row 67
$ws.Range("H67").Value = 45280
$ws.Range("J67").Value = 5452.25
$ws.Range("L67").Value = 5452.25
$ws.Range("N67").Value = -7168.25
# Row 113
$ws.Range("H113").Value = 9912.200000000001
$ws.Range("J113").Value = 7495
$ws.Range("L113").Value = 7495
$ws.Range("N113").Value = -14003
# Row 115
$ws.Range("H115").Value = 1141.5714
$ws.Range("I115").Value = 399.8
$ws.Range("K115").Value = 1199.4
$ws.Range("M115").Value = 367.5999999999999
# Row 125
$ws.Range("H125").Value = 7982.8335
$ws.Range("I125").Value = 17250
$ws.Range("J125").Value = 3349.25
$ws.Range("K125").Value = 155250
$ws.Range("L125").Value = 30143.25
$ws.Range("M125").Value = -152790
$ws.Range("N125").Value = -35063.25
# Row 138
$ws.Range("H138").Value = 4535.4067
$ws.Range("I138").Value = 2654.2727
$ws.Range("K138").Value = 7962.8181
$ws.Range("M138").Value = -2822.8181

$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 6347.797
$ws.Range("I32").Value = 6266.2744
$ws.Range("K32").Value = 6266.2744
$ws.Range("M32").Value = -5979.2744
# Row 45
$ws.Range("H45").Value = 8058.1665
$ws.Range("I45").Value = 8739.799999999999
$ws.Range("K45").Value = 8739.799999999999
$ws.Range("M45").Value = -8362.799999999999
# Row 74
$ws.Range("H74").Value = 3771.625
$ws.Range("I74").Value = 853.63635
$ws.Range("K74").Value = 853.63635
$ws.Range("M74").Value = 20.36365000000001
# Row 77
$ws.Range("H77").Value = 3771.625
$ws.Range("I77").Value = 853.63635
$ws.Range("K77").Value = 4268.18175
$ws.Range("M77").Value = 99.81825000000026

$ws = $wb.Worksheets.Item("BSM")
# Row 22
$ws.Range("H22").Value = 824.5
$ws.Range("I22").Value = 824.5
$ws.Range("K22").Value = 824.5
$ws.Range("M22").Value = -651.5
# Row 74
$ws.Range("H74").Value = 11330.333
$ws.Range("J74").Value = 11330.333
$ws.Range("L74").Value = 11330.333
$ws.Range("N74").Value = -13202.333
# Row 77
$ws.Range("H77").Value = 11330.333
$ws.Range("J77").Value = 11330.333
$ws.Range("L77").Value = 33990.999
$ws.Range("N77").Value = -43350.999
# Row 99
$ws.Range("H99").Value = 28869.5
$ws.Range("I99").Value = 35964.816
$ws.Range("K99").Value = 35964.816
$ws.Range("M99").Value = -34466.816
# Row 141
$ws.Range("H141").Value = 72000
$ws.Range("J141").Value = 72000
$ws.Range("L141").Value = 72000
$ws.Range("N141").Value = -82360

$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 2859.1404
$ws.Range("I31").Value = 2314.25
$ws.Range("K31").Value = 2314.25
$ws.Range("M31").Value = -2019.25
# Row 34
$ws.Range("H34").Value = 2859.1404
$ws.Range("I34").Value = 2314.25
$ws.Range("K34").Value = 2314.25
$ws.Range("M34").Value = -2112.25
# Row 41
$ws.Range("H41").Value = 8784.5
# Row 58
$ws.Range("H58").Value = 4058.8
$ws.Range("J58").Value = 5164.6665
$ws.Range("L58").Value = 5164.6665
$ws.Range("N58").Value = -5570.6665
# Row 62
$ws.Range("H62").Value = 114492.664
$ws.Range("I62").Value = 16553.334
$ws.Range("J62").Value = 212432
$ws.Range("K62").Value = 16553.334
$ws.Range("L62").Value = 212432
$ws.Range("M62").Value = -15929.334
$ws.Range("N62").Value = -213680
# Row 65
$ws.Range("H65").Value = 114492.664
$ws.Range("I65").Value = 16553.334
$ws.Range("J65").Value = 212432
$ws.Range("K65").Value = 82766.67
$ws.Range("L65").Value = 1062160
$ws.Range("M65").Value = -79646.67
$ws.Range("N65").Value = -1068400
# Row 136
$ws.Range("H136").Value = 4058.8
$ws.Range("J136").Value = 5164.6665
$ws.Range("L136").Value = 15493.9995
$ws.Range("N136").Value = -20593.9995

$ws = $wb.Worksheets.Item("CUL")
# Row 3
$ws.Range("H3").Value = 399
$ws.Range("I3").Value = 399
$ws.Range("K3").Value = 1197
$ws.Range("M3").Value = -1085
# Row 16
$ws.Range("H16").Value = 39
$ws.Range("I16").Value = 39
$ws.Range("K16").Value = 117
$ws.Range("M16").Value = 56

$ws = $wb.Worksheets.Item("GSM")
# Row 11
$ws.Range("H11").Value = 7717857
$ws.Range("I11").Value = 7542308
$ws.Range("J11").Value = 9999999
$ws.Range("K11").Value = 7542308
$ws.Range("L11").Value = 9999999
$ws.Range("M11").Value = -7542169
$ws.Range("N11").Value = -10000277
# Row 29
$ws.Range("H29").Value = 13570.429
$ws.Range("J29").Value = 13570.429
$ws.Range("L29").Value = 13570.429
$ws.Range("N29").Value = -14150.429
# Row 35
$ws.Range("H35").Value = 20000
$ws.Range("I35").Value = 20000
$ws.Range("J35").Value = 0
$ws.Range("K35").Value = 20000
$ws.Range("L35").Value = 0
$ws.Range("M35").Value = -19702
$ws.Range("N35").Value = $null
# Row 102
$ws.Range("H102").Value = 6635.0586
$ws.Range("I102").Value = 7843.72
$ws.Range("J102").Value = 3277.6667
$ws.Range("K102").Value = 7843.72
$ws.Range("L102").Value = 3277.6667
$ws.Range("M102").Value = -6221.72
$ws.Range("N102").Value = -6521.6667
# Row 132
$ws.Range("H132").Value = 2869.9688
$ws.Range("I132").Value = 2243.9583
$ws.Range("J132").Value = 4748
$ws.Range("K132").Value = 6731.874899999999
$ws.Range("L132").Value = 14244
$ws.Range("M132").Value = -4201.874899999999
$ws.Range("N132").Value = -19304

$ws = $wb.Worksheets.Item("LTW")
# Row 16
$ws.Range("H16").Value = 2640.2942
$ws.Range("I16").Value = 2682.2856
$ws.Range("K16").Value = 2682.2856
$ws.Range("M16").Value = -2512.2856
# Row 31
$ws.Range("H31").Value = 8334.929
$ws.Range("I31").Value = 2596.9
$ws.Range("J31").Value = 22680
$ws.Range("K31").Value = 2596.9
$ws.Range("L31").Value = 22680
$ws.Range("M31").Value = -2348.9
$ws.Range("N31").Value = -23176
# Row 34
$ws.Range("H34").Value = 5637.6665
$ws.Range("I34").Value = 5208.4
$ws.Range("J34").Value = 6174.25
$ws.Range("K34").Value = 5208.4
$ws.Range("L34").Value = 6174.25
$ws.Range("M34").Value = -5036.4
$ws.Range("N34").Value = -6518.25
# Row 122
$ws.Range("H122").Value = 7159.7676
$ws.Range("I122").Value = 5189.1177
$ws.Range("J122").Value = 14604.444
$ws.Range("K122").Value = 15567.3531
$ws.Range("L122").Value = 43813.33199999999
$ws.Range("M122").Value = -13117.3531
$ws.Range("N122").Value = -48713.33199999999
# Row 132
$ws.Range("H132").Value = 441500.56
$ws.Range("I132").Value = 553630.9399999999
$ws.Range("K132").Value = 1660892.82
$ws.Range("M132").Value = -1658362.82
